$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$tbl = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# Insert one new row inside the Table1 data range at worksheet row 293
# (this pushes the old row 293.. down to 294.. and extends the table
# by one row at the bottom, turning the old "last" row 436 into 437).
# ------------------------------------------------------------------
$ws.Rows.Item(293).Insert()

# The freshly inserted row 293 comes back blank/unstyled; give it the
# same formatting as a normal Table1 data row (copy format only from
# the row directly above it, row 292).
$ws.Range("A292:K292").Copy()
$ws.Range("A293:K293").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the calculated column formula for the new row.
$ws.Range("G293").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# Grow Table1 so it covers the new bottom row too.
$tbl.Resize($ws.Range("A8:K437"))

# Re-evaluate the calculated-column formula on the new last row (437)
# so its cached value isn't left stale/erroring after the resize.
$ws.Range("G437").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# ------------------------------------------------------------------
# Fill in the new leave-card entry data.
# ------------------------------------------------------------------

# Row 292: SL(3-0-0), 1.25 earned, 3 used (W/ Pay), remarks string.
$ws.Range("B292").Value = "SL(3-0-0)"
$ws.Range("C292").Value = 1.25
$ws.Range("H292").Value = 3
$ws.Range("K292").Value = "3/21,23,24/2023"

# Row 293 (brand-new row): SL(1-0-0), 1 used (W/ Pay), dated remark.
$ws.Range("B293").Value = "SL(1-0-0)"
$ws.Range("H293").Value = 1

# K293 needs the date-formatted style (same one already used by K42,
# K47, ... in this sheet) before writing the date serial value, so we
# reuse the existing style instead of growing the style table.
$ws.Range("K42").Copy()
$ws.Range("K293").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K293").Value = 45016

# ------------------------------------------------------------------
# Update the view state to match: scrolled/selected near the new rows.
# ------------------------------------------------------------------
$ws.Range("B294").Select()
